$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths (closest values reachable through the quantized
#     ColumnWidth setter; target stored widths are 29.5703125 /
#     55.85546875 / 68.5703125 characters) ---
$ws.Columns.Item(1).ColumnWidth = 28.736979166666668
$ws.Columns.Item(2).ColumnWidth = 55.022135416666664
$ws.Columns.Item(3).ColumnWidth = 67.73697916666667

# --- New rows of test-report data ---
# testDeletes() sub-tests (rows 10-12)
$ws.Range("A10").Value = "testDeletes() sub-test 1"
$ws.Range("A11").Value = "testDeletes() sub-test 2"
$ws.Range("A12").Value = "testDeletes() sub-test 3"

# testInserts() sub-tests (rows 4-8)
$ws.Range("A4").Value = "testInserts() sub-test 1"
$ws.Range("A5").Value = "testInserts() sub-test 2"
$ws.Range("A6").Value = "testInserts() sub-test 3"
$ws.Range("A7").Value = "testInserts() sub-test 4"

# testUpdates() sub-tests (rows 14-15)
$ws.Range("A14").Value = "testUpdates() sub-test 1"
$ws.Range("A15").Value = "testUpdates() sub-test 2"

# Section header
$ws.Range("A3").Value = "Content Provider (SQLite) tests:"

# Comment column (C) entries
$ws.Range("C7").Value = "Check if all data in the entry is correct."
$ws.Range("C5").Value = """DrinkTable"" should contain only one Drink."
$ws.Range("C4").Value = "Inserts one Drink into the table: ""DrinkTable"" by using the Content Provider."

$ws.Range("A8").Value = "testInserts() sub-test 5"

$ws.Range("C8").Value = "Try to insert a row that already exists. Should return an exception."
$ws.Range("C6").Value = "Moves to the first record with cursor. Should point to the first and only row."
$ws.Range("C10").Value = "Try to delete an empty record."
$ws.Range("C11").Value = "Delete an existing record."
$ws.Range("C12").Value = "Check if the row was actually deleted."
$ws.Range("C14").Value = "Try to update an empty record."
$ws.Range("C15").Value = "Update an existing record."

# Result column (B) entries
$ws.Range("B4").Value = "PASSED! Returned no exceptions."
$ws.Range("B5").Value = "PASSED! Returned one row from table."
$ws.Range("B6").Value = "PASSED! Returned no exceptions."
$ws.Range("B7").Value = "PASSED! Returned correct data."
$ws.Range("B8").Value = "PASSED! Returned an exception."
$ws.Range("B10").Value = "PASSED! Returned zero rows deleted."
$ws.Range("B11").Value = "PASSED! Returned one row deleted."
$ws.Range("B12").Value = "PASSED! Query and check that cursor.getCount returns zero."
$ws.Range("B14").Value = "PASSED! Returned zero rows updated."
$ws.Range("B15").Value = "PASSED! Returned one row updated."

# --- Sheet view / selection ---
$ws.Range("A2").Select()

# --- Page setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
